$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to Text format so numeric-looking values
# (e.g. "0.999", "11.40") are preserved exactly as text, matching the
# original inline-string cell type.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "69.295.34"
$ws.Range("E2").Value = "  +2.47%  "
$ws.Range("D3").Value = "3.384.85"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "588.86"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").Value = "181.38"
$ws.Range("E6").Value = "  +4.18%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.595"
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("D9").Value = "0.197"
$ws.Range("E9").Value = "  +9.55%  "
$ws.Range("D10").Value = "0.588"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").Value = "49.06"
$ws.Range("E11").Value = "  +6.62%  "
$ws.Range("E12").Value = "  +5.04%  "
$ws.Range("D13").Value = "686.88"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").Value = "3.930.21"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").Value = "69.310.64"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").Value = "3.376.90"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("D19").Value = "17.76"
$ws.Range("E19").Value = "  +2.43%  "
$ws.Range("D20").Value = "11.40"
$ws.Range("E20").Value = "  +4.04%  "
$ws.Range("D21").Value = "0.902"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D23").Value = "17.11"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").Value = "105.32"
$ws.Range("E24").Value = "  +6.89%  "
$ws.Range("E25").Value = "  +2.46%  "
$ws.Range("E26").Value = "  +1.80%  "
$ws.Range("E27").Value = "  +2.74%  "
$ws.Range("D28").Value = "34.51"
$ws.Range("E28").Value = "  +3.98%  "
$ws.Range("E29").Value = "  +2.35%  "
$ws.Range("D30").Value = "7.04"
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("B31").Value = "dogwifhat"
$ws.Range("C31").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D31").Value = "3.69"
$ws.Range("E31").Value = "  +11.31%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "11.18"
$ws.Range("E32").Value = "  +2.15%  "
$ws.Range("D33").Value = "557.78"
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("D35").Value = "58.34"
$ws.Range("E35").Value = "  +2.80%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.721.56"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "0.142"
$ws.Range("E38").Value = "  +8.80%  "
$ws.Range("D39").Value = "35.03"
$ws.Range("E39").Value = "  +2.51%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "3.27"
$ws.Range("E40").Value = "  +3.10%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0706"
$ws.Range("E41").Value = "  +5.51%  "
$ws.Range("D42").Value = "2.67"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("E43").Value = "  +1.68%  "
$ws.Range("D44").Value = "0.0418"
$ws.Range("E44").Value = "  +3.48%  "
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("D46").Value = "2.67"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("E48").Value = "  +7.32%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").Value = "132.75"
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("D51").Value = "2.62"
$ws.Range("E51").Value = "  -1.97%  "

"done"
